$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$baseStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "43.328.52"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "2.236.25"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'317.70"
$ws.Range("D5").Style = $baseStyle
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("D6").Value = "'99.03"
$ws.Range("D6").Style = $baseStyle
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "'0.581"
$ws.Range("D7").Style = $baseStyle
$ws.Range("E7").Value = "  +2.54%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").Value = "'37.03"
$ws.Range("D10").Style = $baseStyle
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "'7.69"
$ws.Range("D12").Style = $baseStyle
$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").Value = "'0.865"
$ws.Range("D14").Style = $baseStyle
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").Value = "'14.34"
$ws.Range("D15").Style = $baseStyle
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("D16").Value = "2.231.99"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "43.286.19"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("D18").Value = "'14.21"
$ws.Range("D18").Style = $baseStyle
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("E20").Value = "  +3.00%  "

$ws.Range("D21").Value = "'65.32"
$ws.Range("D21").Style = $baseStyle
$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("E22").Value = "  -2.70%  "

$ws.Range("D23").Value = "'236.25"
$ws.Range("D23").Style = $baseStyle
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("D24").Value = "'2.17"
$ws.Range("D24").Style = $baseStyle
$ws.Range("E24").Value = "  +2.72%  "

$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("E26").Value = "  +3.28%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("D27").Style = $baseStyle
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("E28").Value = "  +2.40%  "

$ws.Range("D29").Value = "'36.76"
$ws.Range("D29").Style = $baseStyle
$ws.Range("E29").Value = "  +10.92%  "

$ws.Range("D30").Value = "'6.37"
$ws.Range("D30").Style = $baseStyle
$ws.Range("E30").Value = "  -2.83%  "

$ws.Range("D31").Value = "'20.29"
$ws.Range("D31").Style = $baseStyle
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").Value = "'0.0871"
$ws.Range("D32").Style = $baseStyle
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").Value = "'157.19"
$ws.Range("D33").Style = $baseStyle
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").Value = "'2.70"
$ws.Range("D34").Style = $baseStyle
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "'3.20"
$ws.Range("D35").Style = $baseStyle
$ws.Range("E35").Value = "  +3.27%  "

$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("E37").Value = "  +1.50%  "

$ws.Range("D38").Value = "'4.39"
$ws.Range("D38").Style = $baseStyle
$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").Value = "'3.69"
$ws.Range("D40").Style = $baseStyle
$ws.Range("E40").Value = "  +3.14%  "

$ws.Range("D41").Value = "'0.0321"
$ws.Range("D41").Style = $baseStyle
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").Value = "'14.33"
$ws.Range("D42").Style = $baseStyle
$ws.Range("E42").Value = "  +19.35%  "

$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").Value = "1.825.61"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("D45").Value = "'0.202"
$ws.Range("D45").Style = $baseStyle
$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").Value = "'83.96"
$ws.Range("D46").Style = $baseStyle
$ws.Range("E46").Value = "  -6.20%  "

$ws.Range("D47").Value = "'5.30"
$ws.Range("D47").Style = $baseStyle
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("D48").Value = "'8.83"
$ws.Range("D48").Style = $baseStyle
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").Value = "'73.87"
$ws.Range("D49").Style = $baseStyle
$ws.Range("E49").Value = "  -4.39%  "

$ws.Range("D50").Value = "'103.23"
$ws.Range("D50").Style = $baseStyle
$ws.Range("E50").Value = "  +0.92%  "

$ws.Range("D51").Value = "'58.10"
$ws.Range("D51").Style = $baseStyle
$ws.Range("E51").Value = "  -4.33%  "
